# Apply scheduled-runner market data refresh to Sheets/Aegis_Profits.xlsx
# (updates currentAveragePrice* / Leve cost / profit columns per sheet)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1085.9114
$ws.Range("I15").Value = 1085.9114
$ws.Range("K15").Value = 3257.7342
$ws.Range("M15").Value = -3088.7342
$ws.Range("H49").Value = 500
$ws.Range("I49").Value = 500
$ws.Range("J49").Value = 500
$ws.Range("K49").Value = 1500
$ws.Range("L49").Value = 1500
$ws.Range("M49").Value = -1364
$ws.Range("N49").Value = -1772
$ws.Range("H103").Value = 5386.5
$ws.Range("J103").Value = 5400.5
$ws.Range("L103").Value = 16201.5
$ws.Range("N103").Value = -17373.5
$ws.Range("H107").Value = 498.1579
$ws.Range("I107").Value = 497.66666
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 497.66666
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1422.33334
$ws.Range("N107").Value = -4340
$ws.Range("H132").Value = 8628187
$ws.Range("I132").Value = 8628187
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 25884561
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -25882031
$ws.Range("H137").Value = 1092.6981
$ws.Range("I137").Value = 1113.4186
$ws.Range("J137").Value = 1003.6
$ws.Range("K137").Value = 3340.2558
$ws.Range("L137").Value = 3010.8
$ws.Range("M137").Value = -790.2557999999999
$ws.Range("N137").Value = -8110.8
$ws.Range("H138").Value = 1784.3422
$ws.Range("I138").Value = 1346.9333
$ws.Range("J138").Value = 3424.625
$ws.Range("K138").Value = 4040.7999
$ws.Range("L138").Value = 10273.875
$ws.Range("M138").Value = 1099.2001
$ws.Range("N138").Value = -20553.875
$ws.Range("H141").Value = 1395.127
$ws.Range("I141").Value = 1230.3833
$ws.Range("J141").Value = 4690
$ws.Range("K141").Value = 3691.1499
$ws.Range("L141").Value = 14070
$ws.Range("M141").Value = 1488.8501
$ws.Range("N141").Value = -24430

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 64017.25
$ws.Range("I2").Value = 1561.909
$ws.Range("J2").Value = 201419
$ws.Range("K2").Value = 1561.909
$ws.Range("L2").Value = 201419
$ws.Range("M2").Value = -1448.909
$ws.Range("N2").Value = -201645
$ws.Range("H10").Value = 1606.6
$ws.Range("I10").Value = 33
$ws.Range("K10").Value = 33
$ws.Range("M10").Value = 137
$ws.Range("H32").Value = 2356.7
$ws.Range("I32").Value = 2296.7742
$ws.Range("J32").Value = 3152.8572
$ws.Range("K32").Value = 2296.7742
$ws.Range("L32").Value = 3152.8572
$ws.Range("M32").Value = -2009.7742
$ws.Range("N32").Value = -3726.8572
$ws.Range("H76").Value = 78978980
$ws.Range("J76").Value = 78978980
$ws.Range("L76").Value = 78978980
$ws.Range("N76").Value = -78979656
$ws.Range("H79").Value = 78978980
$ws.Range("J79").Value = 78978980
$ws.Range("L79").Value = 78978980
$ws.Range("N79").Value = -78981320
$ws.Range("H116").Value = 64017.25
$ws.Range("I116").Value = 1561.909
$ws.Range("J116").Value = 201419
$ws.Range("K116").Value = 1561.909
$ws.Range("L116").Value = 201419
$ws.Range("M116").Value = 732.0909999999999
$ws.Range("N116").Value = -206007

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 64017.25
$ws.Range("I3").Value = 1561.909
$ws.Range("J3").Value = 201419
$ws.Range("K3").Value = 1561.909
$ws.Range("L3").Value = 201419
$ws.Range("M3").Value = -1447.909
$ws.Range("N3").Value = -201647
$ws.Range("H5").Value = 6260
$ws.Range("J5").Value = 6260
$ws.Range("L5").Value = 6260
$ws.Range("N5").Value = -6486
$ws.Range("H20").Value = 26995.281
$ws.Range("I20").Value = 34723.4
$ws.Range("J20").Value = 1234.8889
$ws.Range("K20").Value = 34723.4
$ws.Range("L20").Value = 1234.8889
$ws.Range("M20").Value = -34476.4
$ws.Range("N20").Value = -1728.8889
$ws.Range("H61").Value = 16000
$ws.Range("J61").Value = 16000
$ws.Range("L61").Value = 16000
$ws.Range("N61").Value = -16626

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28483.76
$ws.Range("I31").Value = 1581.0322
$ws.Range("K31").Value = 1581.0322
$ws.Range("M31").Value = -1286.0322
$ws.Range("H34").Value = 28483.76
$ws.Range("I34").Value = 1581.0322
$ws.Range("K34").Value = 1581.0322
$ws.Range("M34").Value = -1379.0322
$ws.Range("H58").Value = 1081.9
$ws.Range("I58").Value = 962.10205
$ws.Range("J58").Value = 1615.5454
$ws.Range("K58").Value = 962.10205
$ws.Range("L58").Value = 1615.5454
$ws.Range("M58").Value = -759.10205
$ws.Range("N58").Value = -2021.5454
$ws.Range("H132").Value = 3089.319
$ws.Range("I132").Value = 2985.1428
$ws.Range("J132").Value = 3393.1667
$ws.Range("K132").Value = 8955.428400000001
$ws.Range("L132").Value = 10179.5001
$ws.Range("M132").Value = -6425.428400000001
$ws.Range("N132").Value = -15239.5001
$ws.Range("H136").Value = 1081.9
$ws.Range("I136").Value = 962.10205
$ws.Range("J136").Value = 1615.5454
$ws.Range("K136").Value = 2886.30615
$ws.Range("L136").Value = 4846.6362
$ws.Range("M136").Value = -336.3061499999999
$ws.Range("N136").Value = -9946.636200000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 65.36364
$ws.Range("I12").Value = 18.666666
$ws.Range("J12").Value = 82.875
$ws.Range("K12").Value = 55.999998
$ws.Range("L12").Value = 248.625
$ws.Range("M12").Value = 117.000002
$ws.Range("N12").Value = -594.625
$ws.Range("H113").Value = 581.931
$ws.Range("J113").Value = 572.86957
$ws.Range("L113").Value = 1718.60871
$ws.Range("N113").Value = -6058.60871
$ws.Range("H122").Value = 490.4
$ws.Range("J122").Value = 387.25
$ws.Range("L122").Value = 3485.25
$ws.Range("N122").Value = -8385.25
$ws.Range("H131").Value = 1301.202
$ws.Range("J131").Value = 1319.5161
$ws.Range("L131").Value = 3958.5483
$ws.Range("N131").Value = -14038.5483

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 700
$ws.Range("I122").Value = 700
$ws.Range("K122").Value = 2100
$ws.Range("M122").Value = 350
$ws.Range("H126").Value = 5167.8335
$ws.Range("J126").Value = 5400
$ws.Range("L126").Value = 16200
$ws.Range("N126").Value = -21140
$ws.Range("H127").Value = 32595
$ws.Range("J127").Value = 32595
$ws.Range("L127").Value = 32595
$ws.Range("N127").Value = -42515
$ws.Range("H132").Value = 2431.98
$ws.Range("I132").Value = 2424.0217
$ws.Range("K132").Value = 7272.0651
$ws.Range("M132").Value = -4742.0651

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2928.5715
$ws.Range("I7").Value = 1412.5
$ws.Range("K7").Value = 1412.5
$ws.Range("M7").Value = -1300.5
$ws.Range("H22").Value = 1092.9231
$ws.Range("J22").Value = 1067.4546
$ws.Range("L22").Value = 1067.4546
$ws.Range("N22").Value = -1657.4546
$ws.Range("H27").Value = 1092.9231
$ws.Range("J27").Value = 1067.4546
$ws.Range("L27").Value = 1067.4546
$ws.Range("N27").Value = -1281.4546
$ws.Range("H46").Value = 920931.6
$ws.Range("I46").Value = 427.6
$ws.Range("J46").Value = 1688018.4
$ws.Range("K46").Value = 427.6
$ws.Range("L46").Value = 1688018.4
$ws.Range("M46").Value = -239.6
$ws.Range("N46").Value = -1688394.4
$ws.Range("H126").Value = 2928.5715
$ws.Range("I126").Value = 1412.5
$ws.Range("K126").Value = 4237.5
$ws.Range("M126").Value = -1767.5
$ws.Range("H132").Value = 4426.8423
$ws.Range("I132").Value = 4483.8887
$ws.Range("K132").Value = 13451.6661
$ws.Range("M132").Value = -10921.6661
$ws.Range("H136").Value = 1322.7106
$ws.Range("I136").Value = 1253.4242
$ws.Range("J136").Value = 1780
$ws.Range("K136").Value = 3760.2726
$ws.Range("L136").Value = 5340
$ws.Range("M136").Value = -1210.2726
$ws.Range("N136").Value = -10440

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1629.5
$ws.Range("I122").Value = 1006.36365
$ws.Range("J122").Value = 2391.111
$ws.Range("K122").Value = 3019.09095
$ws.Range("L122").Value = 7173.333
$ws.Range("M122").Value = -569.0909499999998
$ws.Range("N122").Value = -12073.333
$ws.Range("H132").Value = 1970.0476
$ws.Range("I132").Value = 2196.0852
$ws.Range("J132").Value = 1306.0625
$ws.Range("K132").Value = 6588.2556
$ws.Range("L132").Value = 3918.1875
$ws.Range("M132").Value = -4058.2556
$ws.Range("N132").Value = -8978.1875
$ws.Range("H136").Value = 645
$ws.Range("I136").Value = 405.76315
$ws.Range("J136").Value = 1943.7142
$ws.Range("K136").Value = 1217.28945
$ws.Range("L136").Value = 5831.142599999999
$ws.Range("M136").Value = 1332.71055
$ws.Range("N136").Value = -10931.1426
